# Update "想去人数" (interest count, column F) for the events whose
# counts changed on refresh, across the three sheets that carry this
# column's data: 展览, 演出, 全部类型 (本地生活 / sheet3 is unchanged).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 374
$ws1.Range("F6").Value  = 1239
$ws1.Range("F9").Value  = 184
$ws1.Range("F11").Value = 172
$ws1.Range("F12").Value = 1046
$ws1.Range("F15").Value = 185
$ws1.Range("F16").Value = 1491
$ws1.Range("F19").Value = 346
$ws1.Range("F21").Value = 811
$ws1.Range("F22").Value = 1149
$ws1.Range("F25").Value = 2651
$ws1.Range("F26").Value = 1430
$ws1.Range("F28").Value = 33
$ws1.Range("F29").Value = 392
$ws1.Range("F30").Value = 414
$ws1.Range("F31").Value = 1212
$ws1.Range("F32").Value = 819
$ws1.Range("F33").Value = 1337
$ws1.Range("F34").Value = 157
$ws1.Range("F36").Value = 778
$ws1.Range("F37").Value = 597
$ws1.Range("F38").Value = 666
$ws1.Range("F39").Value = 833
$ws1.Range("F40").Value = 361
$ws1.Range("F41").Value = 243

# Sheet "演出" (Show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 179
$ws2.Range("F15").Value = 622
$ws2.Range("F16").Value = 23
$ws2.Range("F22").Value = 19

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 374
$ws4.Range("F9").Value  = 179
$ws4.Range("F10").Value = 1239
$ws4.Range("F13").Value = 184
$ws4.Range("F16").Value = 172
$ws4.Range("F20").Value = 185
$ws4.Range("F21").Value = 1491
$ws4.Range("F24").Value = 346
$ws4.Range("F26").Value = 1149
$ws4.Range("F27").Value = 2651
$ws4.Range("F29").Value = 1430
$ws4.Range("F32").Value = 33
$ws4.Range("F34").Value = 392
$ws4.Range("F35").Value = 414
$ws4.Range("F36").Value = 1212
$ws4.Range("F39").Value = 819
$ws4.Range("F40").Value = 1337
$ws4.Range("F41").Value = 778
$ws4.Range("F42").Value = 597
$ws4.Range("F43").Value = 666
$ws4.Range("F44").Value = 833
$ws4.Range("F45").Value = 361
$ws4.Range("F46").Value = 19
$ws4.Range("F48").Value = 243
